$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1396103896103896
$ws.Range("C2").Value = 0.6558441558441559
$ws.Range("J2").Value = 0.01298701298701299
$ws.Range("P2").Value = 0.1006493506493507
$ws.Range("S2").Value = 0.09090909090909091
$ws.Range("B3").Value = 0.004807692307692308
$ws.Range("C3").Value = 0.01442307692307692
$ws.Range("J3").Value = 0.02403846153846154
$ws.Range("P3").Value = 0.7596153846153846
$ws.Range("S3").Value = 0.1971153846153846
$ws.Range("B6").Value = 0.05286343612334802
$ws.Range("D6").Value = 0.01762114537444934
$ws.Range("F6").Value = 0.05286343612334802
$ws.Range("J6").Value = 0.2775330396475771
$ws.Range("O6").Value = 0.004405286343612335
$ws.Range("Q6").Value = 0.2202643171806167
$ws.Range("R6").Value = 0.06167400881057269
$ws.Range("S6").Value = 0.3127753303964758
$ws.Range("B7").Value = 0.145
$ws.Range("D7").Value = 0.025
$ws.Range("E7").Value = 0.01
$ws.Range("F7").Value = 0.06
$ws.Range("J7").Value = 0.07000000000000001
$ws.Range("O7").Value = 0.005
$ws.Range("Q7").Value = 0.225
$ws.Range("R7").Value = 0.095
$ws.Range("S7").Value = 0.365
$ws.Range("B8").Value = 0.09328358208955224
$ws.Range("D8").Value = 0.02052238805970149
$ws.Range("F8").Value = 0.07276119402985075
$ws.Range("J8").Value = 0.1175373134328358
$ws.Range("O8").Value = 0.01492537313432836
$ws.Range("Q8").Value = 0.2593283582089552
$ws.Range("R8").Value = 0.07462686567164178
$ws.Range("S8").Value = 0.3470149253731343
$ws.Range("B9").Value = 0.08298755186721991
$ws.Range("D9").Value = 0.01244813278008299
$ws.Range("F9").Value = 0.02489626556016597
$ws.Range("J9").Value = 0.1120331950207469
$ws.Range("O9").Value = 0.01659751037344398
$ws.Range("Q9").Value = 0.2904564315352697
$ws.Range("R9").Value = 0.0912863070539419
$ws.Range("S9").Value = 0.3692946058091287
$ws.Range("B10").Value = 0.1140015302218822
$ws.Range("D10").Value = 0.02142310635042081
$ws.Range("E10").Value = 0.0007651109410864575
$ws.Range("F10").Value = 0.06426931905126243
$ws.Range("J10").Value = 0.09640397857689365
$ws.Range("O10").Value = 0.006885998469778117
$ws.Range("Q10").Value = 0.2930374904361132
$ws.Range("R10").Value = 0.07115531752104055
$ws.Range("S10").Value = 0.3320581484315226
$ws.Range("G11").Value = 0.1506849315068493
$ws.Range("J11").Value = 0.09931506849315068
$ws.Range("K11").Value = 0.2397260273972603
$ws.Range("L11").Value = 0.4931506849315068
$ws.Range("S11").Value = 0.01712328767123288
$ws.Range("G12").Value = 0.7777777777777778
$ws.Range("J12").Value = 0.1437908496732026
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.0457516339869281
$ws.Range("S12").Value = 0.0261437908496732
$ws.Range("G13").Value = 0.7049180327868853
$ws.Range("J13").Value = 0.2950819672131147
$ws.Range("F15").Value = 0.02116402116402116
$ws.Range("H15").Value = 0.1957671957671958
$ws.Range("I15").Value = 0.06878306878306878
$ws.Range("J15").Value = 0.3174603174603174
$ws.Range("K15").Value = 0.09523809523809523
$ws.Range("M15").Value = 0.02116402116402116
$ws.Range("O15").Value = 0.03703703703703703
$ws.Range("S15").Value = 0.2433862433862434
$ws.Range("F16").Value = 0.01851851851851852
$ws.Range("H16").Value = 0.212962962962963
$ws.Range("I16").Value = 0.1018518518518518
$ws.Range("J16").Value = 0.3611111111111111
$ws.Range("K16").Value = 0.07870370370370371
$ws.Range("M16").Value = 0.04166666666666666
$ws.Range("O16").Value = 0.05092592592592592
$ws.Range("S16").Value = 0.1342592592592593
$ws.Range("F17").Value = 0.02932551319648094
$ws.Range("H17").Value = 0.1832844574780059
$ws.Range("I17").Value = 0.1114369501466276
$ws.Range("J17").Value = 0.4281524926686217
$ws.Range("K17").Value = 0.0747800586510264
$ws.Range("M17").Value = 0.02199413489736071
$ws.Range("N17").Value = 0.004398826979472141
$ws.Range("O17").Value = 0.03519061583577713
$ws.Range("S17").Value = 0.1114369501466276
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.1989247311827957
$ws.Range("I18").Value = 0.07526881720430108
$ws.Range("J18").Value = 0.4301075268817204
$ws.Range("K18").Value = 0.1129032258064516
$ws.Range("M18").Value = 0.005376344086021506
$ws.Range("O18").Value = 0.04838709677419355
$ws.Range("S18").Value = 0.1129032258064516
$ws.Range("F19").Value = 0.0170316301703163
$ws.Range("H19").Value = 0.2368207623682076
$ws.Range("I19").Value = 0.09570154095701541
$ws.Range("J19").Value = 0.3584752635847526
$ws.Range("K19").Value = 0.08840227088402271
$ws.Range("M19").Value = 0.0275750202757502
$ws.Range("N19").Value = 0.0008110300081103001
$ws.Range("O19").Value = 0.0689375506893755
$ws.Range("S19").Value = 0.1062449310624493
